$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of the player roster (rows 2-19). Row 1 (header) is unchanged.
$data = @(
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("Brandon Boston Jr.", "SG,SF", "New Orleans Pelicans"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Andre Drummond", "C", "Philadelphia 76ers"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Guerschon Yabusele", "PF,C", "Philadelphia 76ers"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Jordan Poole", "PG,SG", "Washington Wizards")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
